$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2 = 222
$ws.Range("F5").Value2 = 1013
$ws.Range("F6").Value2 = 5667
$ws.Range("F7").Value2 = 525
$ws.Range("F8").Value2 = 739
$ws.Range("F9").Value2 = 975
$ws.Range("F12").Value2 = 40
$ws.Range("F13").Value2 = 596
$ws.Range("F14").Value2 = 35
$ws.Range("F17").Value2 = 1916
$ws.Range("F18").Value2 = 1489
$ws.Range("F19").Value2 = 978
$ws.Range("F21").Value2 = 200
$ws.Range("F22").Value2 = 356
$ws.Range("F23").Value2 = 578
$ws.Range("F24").Value2 = 172
$ws.Range("F25").Value2 = 1062
$ws.Range("F28").Value2 = 3137
$ws.Range("F29").Value2 = 186
$ws.Range("F30").Value2 = 114
$ws.Range("F31").Value2 = 73
$ws.Range("F32").Value2 = 136
$ws.Range("F34").Value2 = 431
$ws.Range("F40").Value2 = 761
$ws.Range("F41").Value2 = 97
$ws.Range("F42").Value2 = 61
$ws.Range("F43").Value2 = 66

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value2 = 220
$ws.Range("F6").Value2 = 153

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value2 = 222
$ws.Range("F5").Value2 = 1013
$ws.Range("F7").Value2 = 5667
$ws.Range("F8").Value2 = 525
$ws.Range("F9").Value2 = 739
$ws.Range("F11").Value2 = 220
$ws.Range("F12").Value2 = 975
$ws.Range("F15").Value2 = 153
$ws.Range("F17").Value2 = 40
$ws.Range("F18").Value2 = 596
$ws.Range("F19").Value2 = 35
$ws.Range("F23").Value2 = 1916
$ws.Range("F24").Value2 = 1489
$ws.Range("F25").Value2 = 978
$ws.Range("F26").Value2 = 200
$ws.Range("F27").Value2 = 356
$ws.Range("F29").Value2 = 578
$ws.Range("F30").Value2 = 172
$ws.Range("F31").Value2 = 1062
$ws.Range("F32").Value2 = 3137
$ws.Range("F33").Value2 = 186
$ws.Range("F34").Value2 = 114
$ws.Range("F35").Value2 = 73
$ws.Range("F36").Value2 = 136
$ws.Range("F38").Value2 = 431
$ws.Range("F43").Value2 = 761
$ws.Range("F44").Value2 = 97
$ws.Range("F45").Value2 = 66
